# Update countries & provincias Spain
# Refreshes the COVID figures for a set of countries (same dataset, later
# timestamp) and re-sorts the table by total cases, descending - exactly
# like the live page did when it re-pulled the feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last refreshed" banner.
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 14:16"

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# 2) Push the refreshed per-country figures (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
#    Row numbers are the countries' current positions before the re-sort.
Set-Row 21  393525 2170 302100 81411 0 15 10014   # Alemania
Set-Row 44  111837 353  97949  12741 0 10 1147    # Oman
Set-Row 46  108969 0    0      0     0 6  5930    # Suecia
Set-Row 55  90380  738  81501  7934  0 4  945     # Bielorrusia
Set-Row 60  71844  2435 53970  16933 0 16 941     # Austria
Set-Row 71  52620  995  29057  22795 0 3  768     # Libia
Set-Row 74  47418  825  40619  6151  0 6  648     # Azerbaiyan
Set-Row 81  37763  760  30877  6192  0 4  694     # Dinamarca
Set-Row 82  37314  999  25989  10274 0 20 1051    # Bosnia y Herzegovina
Set-Row 85  32262  142  27904  3422  0 3  936     # El Salvador
Set-Row 100 17646  1663 7299   10136 0 9  211     # Eslovenia
Set-Row 109 11097  31   10379  414   0 1  304     # Consejo Danes para los Refugiados
Set-Row 124 6028   50   3561   2454  0 0  13      # Sri Lanka
Set-Row 143 4268   38   3098   1159  0 0  11      # Islandia
Set-Row 168 1148   4    1049   64    0 0  35      # Vietnam
Set-Row 182 490    2    473    17    0 0  0       # Islas Feroe
Set-Row 190 282    30   158    123   0 0  1       # Liechtenstein

# 3) Re-sort the whole table (countries only, header/title excluded) by
#    "Casos totales" (column B) descending, the ranking the sheet displays.
$dataRange = $ws.Range("A4:H221")
$keyRange = $ws.Range("B4:B221")
$dataRange.Sort($keyRange, 2)
